$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SOC", "OCV"),
    @(0, 2.5),
    @(0.01, 2.8),
    @(0.05, 3.04),
    @(0.1, 3.14),
    @(0.15, 3.225),
    @(0.2, 3.29),
    @(0.25, 3.35),
    @(0.3, 3.43),
    @(0.35, 3.5),
    @(0.4, 3.56),
    @(0.45, 3.635),
    @(0.5, 3.71),
    @(0.55, 3.77),
    @(0.6, 3.82),
    @(0.65, 3.85),
    @(0.7, 3.89),
    @(0.75, 3.96),
    @(0.8, 4.02),
    @(0.85, 4.06),
    @(0.9, 4.09),
    @(0.95, 4.12),
    @(1, 4.22)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("N17").Select() | Out-Null
